# Weekly update: insert one new week (2 rows: "Primera" and "Segunda" quality
# grades) for Fruta / Piña "Femacal de La Calera" just above the existing
# row 374, pushing the rest of the table down by 2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the current row 374 (old rows 374-394
# become 376-396).
$ws.Rows("374:375").Insert()

# --- New row 374: Primera ---
$ws.Cells.Item(374, 1).Value = 3
$ws.Cells.Item(374, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(374, 3).Value = "Coquimbo"
$ws.Cells.Item(374, 4).Value = 44516
$ws.Cells.Item(374, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(374, 5).Value = 5
$ws.Cells.Item(374, 6).Value = "Fruta"
$ws.Cells.Item(374, 7).Value = 100108
$ws.Cells.Item(374, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(374, 9).Value = 100108005
$ws.Cells.Item(374, 10).Value = "Piña"
$ws.Cells.Item(374, 11).Value = "Caramelo"
$ws.Cells.Item(374, 12).Value = "Primera"
$ws.Cells.Item(374, 13).Value = 108
$ws.Cells.Item(374, 14).Value = 19000
$ws.Cells.Item(374, 15).Value = 19000
$ws.Cells.Item(374, 16).Value = 19000
$ws.Cells.Item(374, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(374, 18).Value = "Ecuador"
$ws.Cells.Item(374, 19).Value = 1583
$ws.Cells.Item(374, 20).Value = 12

# --- New row 375: Segunda ---
$ws.Cells.Item(375, 1).Value = 3
$ws.Cells.Item(375, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(375, 3).Value = "Coquimbo"
$ws.Cells.Item(375, 4).Value = 44516
$ws.Cells.Item(375, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(375, 5).Value = 5
$ws.Cells.Item(375, 6).Value = "Fruta"
$ws.Cells.Item(375, 7).Value = 100108
$ws.Cells.Item(375, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(375, 9).Value = 100108005
$ws.Cells.Item(375, 10).Value = "Piña"
$ws.Cells.Item(375, 11).Value = "Caramelo"
$ws.Cells.Item(375, 12).Value = "Segunda"
$ws.Cells.Item(375, 13).Value = 162
$ws.Cells.Item(375, 14).Value = 19000
$ws.Cells.Item(375, 15).Value = 19000
$ws.Cells.Item(375, 16).Value = 19000
$ws.Cells.Item(375, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(375, 18).Value = "Ecuador"
$ws.Cells.Item(375, 19).Value = 1357
$ws.Cells.Item(375, 20).Value = 14
